# Rebuild the Tool/NamedThing schema sheets:
#  - Tool sheet loses its developer_team/technical_area/keywords/url columns
#    (and their list data-validations) and keeps only id/name/description,
#    matching the NamedThing sheet's shape.
#  - The now-unused ToolCollection sheet (single "entries" header) is dropped.

$wb = $excel.ActiveWorkbook

# Trim the Tool sheet down to id / name / description.
$toolSheet = $wb.Worksheets.Item("Tool")
$toolSheet.Range("A1:G1").ClearContents()
$toolSheet.Range("A1").Value = "id"
$toolSheet.Range("B1").Value = "name"
$toolSheet.Range("C1").Value = "description"

# Remove the data validation drop-downs that applied to the old columns.
$toolSheet.Range("A:G").Validation.Delete()

# Drop the obsolete ToolCollection sheet entirely.
$wb.Worksheets.Item("ToolCollection").Delete()
